$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 data
$ws.Range("A2").Value = 2341760196
$ws.Range("B2").Value = "25-8-2025"
$ws.Range("C2").Value = 0.375

# Add new row 3 data
$ws.Range("A3").Value = 2341760036

$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").Value = "25-8-2025"

$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = 8/24

# Update selection to match target state
$ws.Range("C4").Select() | Out-Null
